# Updated via Streamlit Approval System
# Appends the newly-approved pending-payment record as row 2 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / string columns -------------------------------------------------
$ws.Range("A2").Value  = "WGE 220"
$ws.Range("B2").Value  = "Western Interior Designers & Marine Contractors"
$ws.Range("E2").Value  = "Western Interior Designers & Marine Contractors"
$ws.Range("G2").Value  = "NEFT"
$ws.Range("H2").Value  = "SBIN0003229"
$ws.Range("I2").Value  = "AAAFW8862C"
$ws.Range("J2").Value  = "32AAAFW8862C1Z9"
$ws.Range("K2").Value  = "ANDRIYA THOMAS"
$ws.Range("L2").Value  = "a6499384-2af8-4f12-956a-c1aa1d26a976"
$ws.Range("N2").Value  = "CNRB0000706"
$ws.Range("U2").Value  = "pending"
$ws.Range("X2").Value  = "Payment of stamp paper-2 Nos RPA_UNIQUE_ID : 33531ea5-2d6c-4ab4-896f-eb0746461afe"
$ws.Range("Y2").Value  = "Construction of New RO – Resitement of M/s Jane Austin, Kollam to Puthenkurish, Ernakulam and Conversion from B to A Location: Thiruvaniyoor Village, Kunnathunadu Taluk, Ernakulam District, Puthenkurish – 682308 Divisional Office: Cochin Divisional Office under Kerala State Office"
$ws.Range("Z2").Value  = "PAYMENT OF STAMP PAPER-2 NOS"
$ws.Range("AA2").Value = "executive.westerntender@gmail.com"
$ws.Range("AB2").Value = "ESTIMATION NOT MATCHED"

# --- Numeric columns ---------------------------------------------------------
$ws.Range("D2").Value  = 286962
$ws.Range("F2").Value  = 34413429360
$ws.Range("M2").Value  = 706101053789
$ws.Range("V2").Value  = 1000
$ws.Range("AC2").Value = 0
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0

# --- Date column (C2) --------------------------------------------------------
# Mirrors the original workbook's custom numFmts: id 164 ("yyyy-mm-dd h:mm:ss")
# is registered first (unused by any cell), then id 165 ("YYYY-MM-DD HH:MM:SS")
# is the one actually applied to C2.
$ws.Range("C2").Value = 46297
$ws.Range("C2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("C2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "Row 2 populated"
